# Update Betfair Back/Lay odds on Sheet1 to match the 2026-01-26 refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J2").Value = 3.75
$ws.Range("N2").Value = 2.3
$ws.Range("O2").Value = 1.22
$ws.Range("P2").Value = 2.12
$ws.Range("Q2").Value = 1.68
$ws.Range("R2").Value = 1.44
$ws.Range("S2").Value = 2.58
$ws.Range("U2").Value = 2.46
$ws.Range("F3").Value = 1.52
$ws.Range("G3").Value = 1.65
$ws.Range("H3").Value = 6.6
$ws.Range("I3").Value = 9.4
$ws.Range("J3").Value = 3.95
$ws.Range("K3").Value = 4.7
$ws.Range("L3").Value = 1.01
$ws.Range("M3").Value = 1.01
$ws.Range("N3").Value = 1.01
$ws.Range("O3").Value = 1.32
$ws.Range("P3").Value = 1.78
$ws.Range("R3").Value = 1.1
$ws.Range("S3").Value = 1.01
$ws.Range("T3").Value = 1.01
$ws.Range("U3").Value = 1.01
$ws.Range("V3").Value = 1.11
$ws.Range("W3").Value = 2.52
$ws.Range("X3").Value = 1000
$ws.Range("Y3").Value = 1000
$ws.Range("Z3").Value = 1000
$ws.Range("AA3").Value = 1000
$ws.Range("AB3").Value = 1000
$ws.Range("AC3").Value = 1000
$ws.Range("AD3").Value = 1000
$ws.Range("AE3").Value = 1000
$ws.Range("AF3").Value = 1000
$ws.Range("AG3").Value = 1000
$ws.Range("AH3").Value = 1000
$ws.Range("AI3").Value = 1000
$ws.Range("AJ3").Value = 1000
$ws.Range("AK3").Value = 1000
$ws.Range("AL3").Value = 1000
$ws.Range("AM3").Value = 1000
$ws.Range("AN3").Value = 1000
$ws.Range("AO3").Value = 1000
$ws.Range("I4").Value = 5.4
$ws.Range("J4").Value = 3.25
$ws.Range("K4").Value = 3.8
$ws.Range("L4").Value = 1.01
$ws.Range("M4").Value = 1.01
$ws.Range("N4").Value = 2.26
$ws.Range("O4").Value = 1.45
$ws.Range("R4").Value = 1.19
$ws.Range("S4").Value = 4.2
$ws.Range("T4").Value = 1.01
$ws.Range("U4").Value = 1.01
$ws.Range("V4").Value = 1.23
$ws.Range("W4").Value = 1.86
$ws.Range("X4").Value = 15.5
$ws.Range("Y4").Value = 18.5
$ws.Range("Z4").Value = 50
$ws.Range("AA4").Value = 1000
$ws.Range("AB4").Value = 11
$ws.Range("AC4").Value = 10.5
$ws.Range("AD4").Value = 28
$ws.Range("AE4").Value = 100
$ws.Range("AF4").Value = 17.5
$ws.Range("AG4").Value = 16
$ws.Range("AH4").Value = 30
$ws.Range("AI4").Value = 1000
$ws.Range("AJ4").Value = 38
$ws.Range("AK4").Value = 38
$ws.Range("AL4").Value = 70
$ws.Range("AM4").Value = 1000
$ws.Range("AN4").Value = 30
$ws.Range("AO4").Value = 1000
$ws.Range("G5").Value = 4.6
$ws.Range("I5").Value = 2.08
$ws.Range("K5").Value = 4.4
$ws.Range("L5").Value = 1.38
$ws.Range("Q5").Value = 1.86
$ws.Range("R5").Value = 1.38
$ws.Range("V5").Value = 1.93
$ws.Range("W5").Value = 1.27
$ws.Range("X5").Value = 19
$ws.Range("Y5").Value = 12
$ws.Range("Z5").Value = 15.5
$ws.Range("AA5").Value = 29
$ws.Range("AB5").Value = 19.5
$ws.Range("AC5").Value = 10.5
$ws.Range("AD5").Value = 13
$ws.Range("AE5").Value = 26
$ws.Range("AF5").Value = 40
$ws.Range("AG5").Value = 21
$ws.Range("AH5").Value = 22
$ws.Range("AI5").Value = 44
$ws.Range("AJ5").Value = 110
$ws.Range("AK5").Value = 65
$ws.Range("AL5").Value = 70
$ws.Range("AM5").Value = 120
$ws.Range("AN5").Value = 65
$ws.Range("AO5").Value = 16.5
$ws.Range("I6").Value = 1.46
$ws.Range("J6").Value = 5.2
$ws.Range("F7").Value = 1.18
$ws.Range("G7").Value = 1.22
$ws.Range("H7").Value = 17.5
$ws.Range("I7").Value = 25
$ws.Range("J7").Value = 8.4
$ws.Range("K7").Value = 10.5
$ws.Range("P7").Value = 3.35
$ws.Range("Q7").Value = 1.35
$ws.Range("F8").Value = 1.34
$ws.Range("G8").Value = 1.39
$ws.Range("H8").Value = 8
$ws.Range("W8").Value = 3.55
$ws.Range("G10").Value = 4.4
$ws.Range("H10").Value = 2.22
$ws.Range("I10").Value = 2.42
$ws.Range("J10").Value = 2.92
$ws.Range("K10").Value = 4
$ws.Range("P10").Value = 1.81
$ws.Range("Q10").Value = 1.88
$ws.Range("F12").Value = 1.89
$ws.Range("I12").Value = 4.4
$ws.Range("Q12").Value = 1.32
$ws.Range("J13").Value = 3.7
$ws.Range("L13").Value = 1.24
$ws.Range("R13").Value = 1.74
$ws.Range("S13").Value = 2.08
$ws.Range("U13").Value = 2.8
$ws.Range("G14").Value = 2.04
$ws.Range("I14").Value = 4.8
$ws.Range("J14").Value = 4
$ws.Range("I16").Value = 2.3
$ws.Range("Q17").Value = 1.38
$ws.Range("F21").Value = 1.88
$ws.Range("I21").Value = 5.3
$ws.Range("J21").Value = 3.55
$ws.Range("K21").Value = 3.8
$ws.Range("P21").Value = 1.8
$ws.Range("Q21").Value = 2.06
$ws.Range("P22").Value = 1.62
$ws.Range("Q22").Value = 2.34
$ws.Range("H23").Value = 2.86
$ws.Range("T23").Value = 2.12
$ws.Range("W23").Value = 1.47
$ws.Range("AC23").Value = 6.8
$ws.Range("H24").Value = 3.15
$ws.Range("N25").Value = 3.25
$ws.Range("R25").Value = 1.27
$ws.Range("AH25").Value = 18
$ws.Range("F26").Value = 2.18
$ws.Range("G26").Value = 2.2
$ws.Range("W26").Value = 1.83
$ws.Range("AC26").Value = 7.4
$ws.Range("G27").Value = 1.42
$ws.Range("I27").Value = 12.5
$ws.Range("K27").Value = 5.2
$ws.Range("P27").Value = 1.81
$ws.Range("F28").Value = 2.58
$ws.Range("G28").Value = 2.94
$ws.Range("H28").Value = 2.9
$ws.Range("I28").Value = 3.35
$ws.Range("J28").Value = 2.78
$ws.Range("K28").Value = 3.7
$ws.Range("P28").Value = 1.54
$ws.Range("Q28").Value = 1.01
